# New PO forecast model
# Updates the three report sheets ("Weekly Quantity", "Monthly Trend",
# "PO Forecast") produced by the PO forecasting pipeline for ASIN
# B0C8W4N4TM: one new trailing data point is appended to the weekly and
# monthly aggregates, and the forward-looking forecast sheet is re-pointed
# at the newer model run (later forecast horizon + revised quantities).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append the newest observed week
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A21").NumberFormat = $wsWeekly.Range("A20").NumberFormat
$wsWeekly.Range("A21").Value = 45662.99999999999
$wsWeekly.Range("B21").Value = 2

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append the newest observed month
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A12").NumberFormat = $wsMonthly.Range("A11").NumberFormat
$wsMonthly.Range("A12").Value = 45688.99999999999
$wsMonthly.Range("B12").Value = 2

# ---------------------------------------------------------------------
# Sheet 3: "PO Forecast" - refreshed forecast model output
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Revised forecast quantities for already-known historical weeks
$wsForecast.Range("B3").Value = 15
$wsForecast.Range("B10").Value = 7
$wsForecast.Range("B14").Value = 6
$wsForecast.Range("B15").Value = 6
$wsForecast.Range("B16").Value = 4
$wsForecast.Range("B18").Value = 3

# The forward forecast horizon (rows 21-28) rolls forward by three weeks;
# shift each existing forecast date and append one new trailing week (29).
$wsForecast.Range("A21").Value = 45662.99999999999
$wsForecast.Range("A22").Value = 45669.99999999999
$wsForecast.Range("A23").Value = 45676.99999999999
$wsForecast.Range("A24").Value = 45683.99999999999
$wsForecast.Range("A25").Value = 45690.99999999999
$wsForecast.Range("A26").Value = 45697.99999999999
$wsForecast.Range("A27").Value = 45704.99999999999
$wsForecast.Range("A28").Value = 45711.99999999999

$wsForecast.Range("A29").NumberFormat = $wsForecast.Range("A28").NumberFormat
$wsForecast.Range("A29").Value = 45718.99999999999
$wsForecast.Range("B29").Value = 0
